$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete obsolete rows 6-8 (data now only spans rows 1-5)
$ws.Rows("6:8").Delete()

# Column A holds dates stored as text; force text format so Excel does not
# auto-convert the strings into date serials.
$ws.Range("A2:A5").NumberFormat = "@"

$ws.Range("A2").Value = "2025-08-30"
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 88
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 880
$ws.Range("G2").Value = 88
$ws.Range("H2").Value = 88
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 89
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 104
$ws.Range("M2").Value = 184
$ws.Range("N2").Value = 96
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 200
$ws.Range("A3").Value = "2025-08-31"
$ws.Range("B3").Value = 90
$ws.Range("C3").Value = 800
$ws.Range("D3").Value = 70
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 80
$ws.Range("G3").Value = 99
$ws.Range("H3").Value = 76
$ws.Range("I3").Value = 78
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 98
$ws.Range("L3").Value = 240
$ws.Range("M3").Value = 239
$ws.Range("N3").Value = 6865
$ws.Range("O3").Value = 176
$ws.Range("P3").Value = 7105
$ws.Range("Q3").Value = 415
$ws.Range("A4").Value = "2025-07-01"
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 124529
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 27
$ws.Range("M4").Value = 108
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 18
$ws.Range("P4").Value = 36
$ws.Range("Q4").Value = 126
$ws.Range("A5").Value = "2025-09-01"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 9
$ws.Range("H5").Value = 9
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 9
$ws.Range("L5").Value = 27
$ws.Range("M5").Value = 27
$ws.Range("N5").Value = 18
$ws.Range("O5").Value = 18
$ws.Range("P5").Value = 45
$ws.Range("Q5").Value = 45
